# "Generate Report for Archive"
# The localization status report regenerated: the "Status" cells that used
# to read "Ready for handoff" now read "In Translation" (same shared text
# is used on the Overview summary sheet and on each per-locale sheet), and
# the now-narrower text lets the "Status" columns shrink to fit.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status text (shared by the summary sheet's two locale columns
# and by the Status column on each per-locale detail sheet).
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Re-fit the "Status" columns now that the text is shorter.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
